$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new blank row above row 9; existing rows 9..54 shift down to 10..55
$ws.Rows(9).Insert()

# Give the new row 9 the same cell formatting as the row directly below it
# (the original row 9 content, which has shifted down to row 10).
$ws.Range("A10:G10").Copy()
$ws.Range("A9:G9").PasteSpecial(-4122)  # xlPasteFormats
$excel.CutCopyMode = 0

# Column A has no entry on this milestone row
$ws.Range("A9").Clear()

# The row that used to be row 9 (now row 10) had a stray one-off style on A10
# that duplicated the plain "vertical top + wrap" style used elsewhere (A11).
# Normalise it so the now-unused duplicate style can drop out of the stylesheet.
$ws.Range("A11").Copy()
$ws.Range("A10").PasteSpecial(-4122)  # xlPasteFormats
$excel.CutCopyMode = 0

# New milestone entry content
$ws.Range("B9").Value = "Not a milestone - Implemented functions for database.c. Also made changes to avl_tree to accommodate the database better. linked_list had a function name refactored as well for clarity"
$ws.Range("C9").Value = 45774
$ws.Range("D9").Value = "insert_list`nprint_avl_node_without_tree`nprint_avl_without_tree`nprint_events`nprint_fighters_in_event`nadd_fighter_list`nadd_fighter_avl`nadd_fighter_to_event`nadd_event`nprint_events)of_fighter_list`nprint_events_of_fighter_avl`nprint_events_of_fighter`ntest_database"
$ws.Range("E9").Value = "No other resources were used for this commit."
$ws.Range("F9").Value = "test_database - This currently only includes unstructured adhoc tests that aimed to prove that the functions work in their most basic use cases. Functions in database.c require further testing"
$ws.Range("G9").Value = "It has become obvious across multiple files that the commenting style is somewhat inconsistent. There is (I feel) a good amount of comments written, but the style of them and the display needs to be made uniform. There are also some functions that have far less commenting than others. A commit will need to focus solely on updating comments.`nIt has also come to my attention that my use of a doubly linked list may have been somewhat unnecessary. It was meant to cut down the time it takes to search for records. But later on I realised that it would be pointless when searching for fighters, as the alphabetical order of names can be quite random. I found that it would only be useful for searching from the tail or head of the list that contains event codes, as if the tail was an event such as 315, and the user was searching for event 310, then it can be easily calculated that it would be better to go from the tail. However the most time complex function of the database, which is finding which events a fighter has fought in, requires every event to be searched. This means that the usefulness of the doubly linked list is lessened."

# Custom row height to fit the long note
$ws.Rows(9).RowHeight = 246.75

# Restore the view: scroll back to the top-left and select the originally
# active cell, which shifted down from G9 to G10 along with the data.
$ws.Application.ActiveWindow.ScrollRow = 1
$ws.Application.ActiveWindow.ScrollColumn = 1
$ws.Range("G10").Select()

Write-Host "Done"
